$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (run_time) values for rows 2-11
$cValues = @(0.4687700271606445, 0.41029953956604, 0.3593680858612061, 0.3749837875366211, 0.375016450881958, 0.3749837875366211, 0.3593878746032715, 0.374950647354126, 0.3750007152557373, 0.3593332767486572)
for ($i = 0; $i -lt $cValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $cValues[$i]
}

# Update columns E:Y (max_er and iter 0..19) values for rows 2-11
$data = @(
    @(723.12269821462, 0.02895402551606931, 0.02420724131049888, 0.02254469845991823, 0.02021272996409079, 0.0178570868321924, 0.01692457362821394, 0.01634042763285096, 0.01558544942255055, 0.0152741554346504, 0.01489683376842072, 0.01473251643351146, 0.01450267481986731, 0.01442613643453485, 0.01431090448550254, 0.01426087019712422, 0.01418276438016629, 0.01415993859647826, 0.01412060744088734, 0.01410938405110328, 0.01409595902952475),
    @(734.7025645082249, 0.02830746954856514, 0.0241847173594625, 0.02132354969162452, 0.01891710893446091, 0.01775987892090011, 0.01710411123993955, 0.01648082843693811, 0.01614211835279102, 0.01551180436105457, 0.01545202371889318, 0.01498945156810305, 0.01483560466547924, 0.01473618557322881, 0.01458641858780619, 0.01449243630502437, 0.01441644410916628, 0.014399366135601, 0.01436891643831896, 0.01432168741731432, 0.01432168741731432),
    @(728.1761857537495, 0.02792232367933633, 0.02412197968177879, 0.02190576185587806, 0.0199845257569125, 0.01890970805787141, 0.01734734060867298, 0.01668550292902765, 0.01628836590767601, 0.01572714546728058, 0.01518862304194023, 0.01485212261573434, 0.01469996860253451, 0.01461310907172879, 0.01448740473195386, 0.01441322318888309, 0.01433306866601924, 0.0142636602691798, 0.0142253718735334, 0.01420575579999149, 0.01419446755855262),
    @(741.9892094092465, 0.02856645159225074, 0.02375441974663346, 0.02159800308042014, 0.01962318447787751, 0.01854096582538122, 0.01755848685089207, 0.01651517518550831, 0.01630370439549542, 0.01585014624891508, 0.01542934049596434, 0.01514152384943673, 0.01489030518777658, 0.01479669607490839, 0.01472516596879866, 0.01463181419086085, 0.01458606510925597, 0.01452891759038102, 0.01450106574075306, 0.01448298807965358, 0.01446372727893268),
    @(730.4205841133407, 0.02820470132401782, 0.02385439514809048, 0.02120874465893905, 0.01883387010721196, 0.01828127489755379, 0.01688192084886573, 0.01623738746065336, 0.0157297890085198, 0.01537786681772963, 0.01510936356747029, 0.01492175946943699, 0.01469065447790917, 0.01455735128326695, 0.01451743438242823, 0.01442915495829608, 0.01434659947891437, 0.01431591871172111, 0.01428221362550966, 0.01425574607857233, 0.01423821801390527),
    @(717.3933043366033, 0.02895424877642565, 0.0243450404970844, 0.0208812413620671, 0.01827512795862534, 0.01719364834559635, 0.01656543315311258, 0.01551821134968519, 0.01495214729861681, 0.01486765969788954, 0.01453847172487082, 0.01449131851335604, 0.01429024316660234, 0.01424574550991967, 0.01418005601955564, 0.01413364906748618, 0.01409173395000661, 0.01404148561443623, 0.01401804516691865, 0.01398911320186583, 0.01398427493833534),
    @(736.8859993974511, 0.02844094105898348, 0.0238872106548345, 0.02117303240905511, 0.01976317855140996, 0.0174184318678946, 0.01704819842531934, 0.01607281020953092, 0.01554818627074494, 0.01525393762227005, 0.0150493070154416, 0.01494346579501621, 0.0148249433600932, 0.01470935863378267, 0.01463988200628815, 0.01450019721960969, 0.01444775592107629, 0.01441908575301075, 0.01439117295173016, 0.0143797357749662, 0.01436424950092497),
    @(743.8537593683468, 0.02898418242965929, 0.02480256062002152, 0.02175066973550077, 0.01972671083275803, 0.01867787567974667, 0.01776363103671519, 0.01716687673894971, 0.01668535684981325, 0.01594994701705937, 0.01573356148984531, 0.01537766266455609, 0.01514868651983574, 0.01479562340121629, 0.01476935435573229, 0.01469065510022463, 0.01459509747326514, 0.01456438878478717, 0.01454886160017186, 0.01451366283134936, 0.01450007328203405),
    @(740.1754453887843, 0.0275811010245802, 0.02336871173321707, 0.02095798671584988, 0.01920414691266778, 0.01798324948101347, 0.01732452380147715, 0.0167004243629385, 0.01621349358808219, 0.01583294966137693, 0.01546990651578409, 0.01510601107858516, 0.01506826300003977, 0.0149175429955575, 0.01476520690423834, 0.01469252033085005, 0.01457010111342898, 0.0145264320425055, 0.01446210386863472, 0.01442837125514199, 0.01442837125514199),
    @(732.147887553183, 0.02939179100214321, 0.02453191134325624, 0.02111672236015432, 0.01931652951095221, 0.01782319697706157, 0.01672643206869538, 0.01622414793962052, 0.01531166190089913, 0.01531166190089913, 0.01503654114710405, 0.01484640139175698, 0.01470033301702791, 0.01461966576290391, 0.01450653360183916, 0.01442785006591864, 0.01440446068207449, 0.01433453826006739, 0.0143147745045003, 0.01428174495130035, 0.01427188864626087)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 5).Value = $data[$r][$c]
    }
}
